# Fix: Remove duplicate 'Withdrawal Amount' column from AEPS template.
# The table has 6 columns, with the 5th ("Withdrawal Date") and 6th
# ("Withdrawal Amount") columns mislabeled/duplicated - the header row's
# 5th and 6th cells both read "Withdrawal Amount" while the correct
# header sequence should be: Sl. No. | Account No. | Transaction ID |
# Withdrawal Date | Withdrawal Amount. Delete the extra 6th column
# (and its data in each row) to restore the correct 5-column table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Delete the 6th (last, extra) column across every row of the table.
$t.Columns.Item(6).Delete()
